# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" worksheet right after "总计" (and before "2022-Q2"),
#    populated with the quarter's fund-holdings table.
# 2) Insert a new row at the top of the "总计" (summary) sheet's data with the
#    2022-Q3 totals, shifting every existing quarter row down by one and
#    renumbering the running index in column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: new "2022-Q3" worksheet
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$q2Sheet      = $wb.Worksheets.Item("2022-Q2")

$q3Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q3Sheet.Name = "2022-Q3"

# Borrow the header / index-column formatting from the 2022-Q2 sheet so the
# new sheet matches the workbook's existing look (bold, centered, bordered).
$q2Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2:A6").Copy()
$q3Sheet.Range("A2:A6").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# index, fund code, fund name, fund size, stock position, position pct,
# held market value (亿元), position rank
$q3Data = @(
    @(0, "'162006", "长城久富核心成长混合（LOF）A", "'31.61", "'71.49", "'3.67", "'1.1601", 6),
    @(1, "'015383", "长城久富核心成长混合（LOF）C", "'1.82",  "'71.49", "'3.67", "'0.0668", 6),
    @(2, "'000976", "长城新兴产业灵活配置混合",     "'1.53",  "'71.73", "'3.68", "'0.0563", 6),
    @(3, "'001707", "诺安高端制造股票A",            "'1.12",  "'90.57", "'3.47", "'0.0389", 10),
    @(4, "'014536", "诺安高端制造股票C",            "'0.00",  "'90.57", "'3.47", 0,         10)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]
    $q3Sheet.Cells.Item($r, 1).Value = $row[0]
    $q3Sheet.Cells.Item($r, 2).Value = $row[1]
    $q3Sheet.Cells.Item($r, 3).Value = $row[2]
    $q3Sheet.Cells.Item($r, 4).Value = $row[3]
    $q3Sheet.Cells.Item($r, 5).Value = $row[4]
    $q3Sheet.Cells.Item($r, 6).Value = $row[5]
    $q3Sheet.Cells.Item($r, 7).Value = $row[6]
    $q3Sheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# Part 2: update the "总计" summary sheet
# ---------------------------------------------------------------------------
$ws = $summarySheet

# Push the existing data rows (2..8) down to (3..9) and insert a fresh row 2.
$ws.Range("A2").EntireRow.Insert()
$ws.Range("B2:D2").ClearFormats()

# Give the new A2 the same "index column" look (style) as the rest of column A.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 1.32

# Renumber the running index in column A for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
